$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the steel description in B2: remove the "/RME" segment from the
# "16% S/LFM+CDH/RME/H:1" line (commit: "simplify steel description (remove RME)").
$current = $ws.Range("B2").Value2
$updated = $current.Replace("16% S/LFM+CDH/RME/H:1", "16% S/LFM+CDH/H:1")
$ws.Range("B2").Value = $updated

# The long multi-line description now wraps within the cell.
$ws.Range("B2").WrapText = $true

# Let the row grow to fit the wrapped text (it hits Excel's row-height cap).
$ws.Rows.Item(2).RowHeight = 409.6

# The user's selection ends up spanning B2:B13.
$ws.Range("B2:B13").Select()
